$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Update PRO number
$ws.Range("B3").Value = 20578803

# Update Instrument SN / ICA SN values (these feed Template_printout via formulas,
# and B7's MOD formula recalculates automatically off the new B4 value)
$ws.Range("B4").Value = "A01603"
$ws.Range("B5").Value = "APXCAS2131011"

# Match the author's recorded UI selection on the Input sheet, then
# restore the originally-active sheet (Template_printout) so the
# workbook's active-tab state is unchanged
$ws.Range("B5").Select() | Out-Null
$wb.Worksheets.Item("Template_printout").Activate() | Out-Null
